$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111 (shifts existing rows 111..210 down to 112..211)
$ws.Rows.Item(111).Insert()

# Fill the new row 111 with its data (same template columns as the
# surrounding Jengibre / Vega Modelo de Temuco rows, new unique values)
$ws.Range("A111").Value = 10
$ws.Range("B111").Value = "Vega Modelo de Temuco"
$ws.Range("C111").Value = "La Araucanía"
$ws.Range("D111").Value = 44880
$ws.Range("E111").Value = 9
$ws.Range("F111").Value = 100114007
$ws.Range("G111").Value = "Jengibre"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 40
$ws.Range("K111").Value = 18000
$ws.Range("L111").Value = 20000
$ws.Range("M111").Value = 19000
$ws.Range("N111").Value = "$/caja 13 kilos"
$ws.Range("O111").Value = "Perú"
$ws.Range("P111").Value = 1462
$ws.Range("Q111").Value = 13
$ws.Range("R111").Value = "Hortaliza"
